$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New "Math" section header in C21 (plain text, matches the look of the
# other section headers such as C18 "Strings").
$c21 = $ws1.Range("C21")
$c21.Value = "Math"
$ws1.Rows("21:21").RowHeight = 15.75

# New hyperlink cell C22 pointing at the palindrome-number problem.
$c22 = $ws1.Range("C22")
$c22.Value = "https://leetcode.com/problems/palindrome-number"
$ws1.Rows("22:22").RowHeight = 15.75
$ws1.Hyperlinks.Add($c22, "https://leetcode.com/problems/palindrome-number") | Out-Null
$c22.Style = "Hyperlink"

# Add the new, empty "Sheet2" right after "Sheet1".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Restore Sheet1 as the active/selected sheet (Worksheets.Add activates
# the freshly inserted sheet as a side effect).
$ws1.Activate()
